# Update "想去人数" (want-to-go count) figures in the 展览 (sheet1) and
# 全部类型 (sheet4) sheets, reflecting a newer scrape of the source data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet (rows keyed by their current row number) ---
$wsExhibit.Range("F2").Value  = 13892   # was 13863
$wsExhibit.Range("F6").Value  = 496     # was 497
$wsExhibit.Range("F7").Value  = 1212    # was 1210
$wsExhibit.Range("F8").Value  = 1036    # was 1035
$wsExhibit.Range("F9").Value  = 13907   # was 13906
$wsExhibit.Range("F10").Value = 14836   # was 14816
$wsExhibit.Range("F11").Value = 2       # was 1
$wsExhibit.Range("F20").Value = 22      # was 21
$wsExhibit.Range("F23").Value = 1151    # was 1149
$wsExhibit.Range("F26").Value = 5763    # was 5757
$wsExhibit.Range("F28").Value = 1062    # was 1061
$wsExhibit.Range("F29").Value = 5425    # was 5422
$wsExhibit.Range("F30").Value = 50      # was 49
$wsExhibit.Range("F32").Value = 287     # was 281

# --- 全部类型 sheet (same events, offset by one extra row) ---
$wsAll.Range("F2").Value  = 13892   # was 13863
$wsAll.Range("F7").Value  = 496     # was 497
$wsAll.Range("F8").Value  = 1212    # was 1210
$wsAll.Range("F9").Value  = 1036    # was 1035
$wsAll.Range("F10").Value = 13907   # was 13906
$wsAll.Range("F11").Value = 14836   # was 14816
$wsAll.Range("F12").Value = 2       # was 1
$wsAll.Range("F21").Value = 22      # was 21
$wsAll.Range("F24").Value = 1151    # was 1149
$wsAll.Range("F27").Value = 5763    # was 5757
$wsAll.Range("F29").Value = 1062    # was 1061
$wsAll.Range("F30").Value = 5425    # was 5422
$wsAll.Range("F31").Value = 50      # was 49
$wsAll.Range("F33").Value = 287     # was 281
